# The source workbook has two sheets:
#   - "produit (2)" (sheetId 2): holds a Power Query-backed table/ListObject
#     ("produit") whose data was pulled from an external query connection.
#   - "produit" (sheetId 1): a plain sheet with just two text cells (A1/A2)
#     holding the header row and (one) data row of the scraped CSV, pasted
#     as raw text.
#
# The commit ("simplification de l'ecriture + ecriture du csv") drops the
# Power-Query-driven sheet/table entirely and keeps only the simple
# "produit" sheet, refreshing its two cells with a freshly exported CSV
# (different column order/names, real UPC/url columns, etc).

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# Remove the query/table sheet "produit (2)" completely.
$wb.Worksheets.Item("produit (2)").Delete()

# "produit" is now the only worksheet; make sure it is the active/selected
# one (it becomes rId1 / first and only tab).
$ws = $wb.Worksheets.Item("produit")
$ws.Select()

# Refresh the two text cells with the current CSV export: header line, then
# the single data row (re-ordered/renamed columns vs. the previous export).
$ws.Range("A1").Value = 'title,price_including_tax,product_description,number_available,category,price_excluding_tax,UPC,product_url,image_url,review_rating'
$ws.Range("A2").Value = 'http://books.toscrape.com/catalogue/a-light-in-the-attic_1000/index.html,http://books.toscrape.com/media/cache/fe/72/fe72f0532301ec28892ae79a629a293c.jpg,A Light in the Attic,In stock (22 available),Poetry,a897fe39b1053632,"It''s hard to imagine a world without A Light in the Attic. This now-classic collection of poetry and drawings from Shel Silverstein celebrates its 20th anniversary with this special edition. Silverstein''s humorous and creative verse can amuse the dowdiest of readers. Lemon-faced adults and fidgety kids sit still and read these rhythmic words and laugh and smile and love th It''s hard to imagine a world without A Light in the Attic. This now-classic collection of poetry and drawings from Shel Silverstein celebrates its 20th anniversary with this special edition. Silverstein''s humorous and creative verse can amuse the dowdiest of readers. Lemon-faced adults and fidgety kids sit still and read these rhythmic words and laugh and smile and love that Silverstein. Need proof of his genius? RockabyeRockabye baby, in the treetopDon''t you know a treetopIs no safe place to rock?And who put you up there,And your cradle, too?Baby, I think someone down here''sGot it in for you. Shel, you never sounded so good. ...more",Ã‚Â£51.77,Threeout of five'

# Match the saved selection state (both cells selected).
$ws.Range("A1:A2").Select()
